$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 538, shifting existing rows 538:629 down to 539:630.
$ws.Rows.Item(538).Insert()

# Populate the newly inserted row 538 with the new record.
$ws.Range("A538").Value = 5
$ws.Range("B538").Value = "Macroferia Regional de Talca"
$ws.Range("C538").Value = "Maule"
$ws.Range("D538").Value = 45180
$ws.Range("E538").Value = 7
$ws.Range("F538").Value = 100112023
$ws.Range("G538").Value = "Brócoli"
$ws.Range("H538").Value = "Sin especificar"
$ws.Range("I538").Value = "Primera"
$ws.Range("J538").Value = 3000
$ws.Range("K538").Value = 800
$ws.Range("L538").Value = 800
$ws.Range("M538").Value = 800
$ws.Range("N538").Value = "`$/unidad"
$ws.Range("O538").Value = "Región del Maule"
$ws.Range("P538").Value = 800
$ws.Range("Q538").Value = 1
$ws.Range("R538").Value = "Hortaliza"
